# Add 2022-Q3 data
#
# 1. Insert a brand-new worksheet named "2022-Q3" right after the "总计" sheet,
#    holding the same kind of per-fund breakdown table as the other quarter sheets.
# 2. Update the "总计" (summary) sheet: insert a new top data row for 2022-Q3 and
#    push the existing quarters down, appending a duplicated 2021-Q1 row at the end.

$wb = $excel.ActiveWorkbook

$totalSheet = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1) Create the new "2022-Q3" sheet right after "总计", using the "2022-Q2"
#    sheet (same column layout / header / styling) as a formatting template.
# ---------------------------------------------------------------------------
$q3 = $wb.Worksheets.Add($null, $totalSheet)
$q3.Name = "2022-Q3"

# NOTE: the lookup of the template sheet must happen *after* the new sheet
# has been inserted into the workbook, otherwise the by-name reference can
# end up stale/pointing at the wrong worksheet.
$templateSheet = $wb.Worksheets.Item("2022-Q2")

$templateSheet.Range("A1:H7").Copy()
$q3.Range("A1").PasteSpecial(-4122)  # xlPasteFormats (styling only, no text/values)
$q3.Range("A1").Clear()  # the template has no A1 cell at all (header starts at B1)

# Header row (same column headers as the other quarter sheets)
$q3.Range("B1").Value = "基金代码"
$q3.Range("C1").Value = "基金名称"
$q3.Range("D1").Value = "基金规模"
$q3.Range("E1").Value = "股票总仓位"
$q3.Range("F1").Value = "仓位占比"
$q3.Range("G1").Value = "持有市值(亿元)"
$q3.Range("H1").Value = "仓位排名"

# Make sure the value-like text columns (B, D, E, F, G) are stored as plain text,
# matching how the source data is represented elsewhere in the workbook.
$q3.Range("B2:B7").NumberFormat = "@"
$q3.Range("D2:G7").NumberFormat = "@"

# Row 2
$q3.Range("A2").Value = 0
$q3.Range("B2").Value = "001487"
$q3.Range("C2").Value = "宝盈优势产业灵活配置混合A"
$q3.Range("D2").Value = "10.11"
$q3.Range("E2").Value = "91.85"
$q3.Range("F2").Value = "2.99"
$q3.Range("G2").Value = "0.3023"
$q3.Range("H2").Value = 5

# Row 3
$q3.Range("A3").Value = 1
$q3.Range("B3").Value = "013895"
$q3.Range("C3").Value = "宝盈成长精选混合A"
$q3.Range("D3").Value = "8.51"
$q3.Range("E3").Value = "90.34"
$q3.Range("F3").Value = "2.83"
$q3.Range("G3").Value = "0.2408"
$q3.Range("H3").Value = 10

# Row 4
$q3.Range("A4").Value = 2
$q3.Range("B4").Value = "001075"
$q3.Range("C4").Value = "宝盈转型动力灵活配置混合A"
$q3.Range("D4").Value = "4.35"
$q3.Range("E4").Value = "91.90"
$q3.Range("F4").Value = "3.14"
$q3.Range("G4").Value = "0.1366"
$q3.Range("H4").Value = 10

# Row 5
$q3.Range("A5").Value = 3
$q3.Range("B5").Value = "012771"
$q3.Range("C5").Value = "宝盈优势产业灵活配置混合C"
$q3.Range("D5").Value = "3.62"
$q3.Range("E5").Value = "91.85"
$q3.Range("F5").Value = "2.99"
$q3.Range("G5").Value = "0.1082"
$q3.Range("H5").Value = 5

# Row 6
$q3.Range("A6").Value = 4
$q3.Range("B6").Value = "013896"
$q3.Range("C6").Value = "宝盈成长精选混合C"
$q3.Range("D6").Value = "2.93"
$q3.Range("E6").Value = "90.34"
$q3.Range("F6").Value = "2.83"
$q3.Range("G6").Value = "0.0829"
$q3.Range("H6").Value = 10

# Row 7
$q3.Range("A7").Value = 5
$q3.Range("B7").Value = "015389"
$q3.Range("C7").Value = "宝盈转型动力灵活配置混合C"
$q3.Range("D7").Value = "0.24"
$q3.Range("E7").Value = "91.90"
$q3.Range("F7").Value = "3.14"
$q3.Range("G7").Value = "0.0075"
$q3.Range("H7").Value = 10

# ---------------------------------------------------------------------------
# 2) Update the "总计" sheet: shift existing quarter rows down by one and add
#    the new 2022-Q3 row on top, plus a new trailing 2021-Q1 row.
# ---------------------------------------------------------------------------

# Duplicate the formatting of the last existing data row (row 6) into the new
# row 7 before we touch any values, so the new row picks up the same styling
# (bold index column etc.) as the rest of the table.
$totalSheet.Range("A6:D6").Copy()
$totalSheet.Range("A7:D7").PasteSpecial(-4122)

# New bottom row: 2021-Q1 (same figures as the old 2021-Q2 row, 2 / 0.08)
$totalSheet.Range("A7").Value = 5
$totalSheet.Range("B7").Value = "2021-Q1"
$totalSheet.Range("C7").Value = 2
$totalSheet.Range("D7").Value = 0.08

# Row 6: 2021-Q2 (was row 5)
$totalSheet.Range("A6").Value = 4
$totalSheet.Range("B6").Value = "2021-Q2"
$totalSheet.Range("C6").Value = 2
$totalSheet.Range("D6").Value = 0.08

# Row 5: 2021-Q4 (was row 4)
$totalSheet.Range("A5").Value = 3
$totalSheet.Range("B5").Value = "2021-Q4"
$totalSheet.Range("C5").Value = 6
$totalSheet.Range("D5").Value = 0.43

# Row 4: 2022-Q1 (was row 3)
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("B4").Value = "2022-Q1"
$totalSheet.Range("C4").Value = 6
$totalSheet.Range("D4").Value = 0.42

# Row 3: 2022-Q2 (was row 2)
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("B3").Value = "2022-Q2"
$totalSheet.Range("C3").Value = 15
$totalSheet.Range("D3").Value = 1.94

# Row 2: new 2022-Q3 row
$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q3"
$totalSheet.Range("C2").Value = 6
$totalSheet.Range("D2").Value = 0.88
